$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()

$ws.Range("H33").Value = 50174.9
$ws.Range("I33").Value = 50174.9
$ws.Range("K33").Value = 50174.9
$ws.Range("M33").Value = -49945.9

$ws.Range("H107").Value = 400298.97
$ws.Range("I107").Value = 526564.2
$ws.Range("J107").Value = 459.16666
$ws.Range("K107").Value = 526564.2
$ws.Range("L107").Value = 459.16666
$ws.Range("M107").Value = -524644.2
$ws.Range("N107").Value = -4299.16666

$ws.Range("H121").Value = 690.1852
$ws.Range("J121").Value = 681.4
$ws.Range("L121").Value = 2044.2
$ws.Range("N121").Value = -5538.2

$ws.Range("H141").Value = 1771.909
$ws.Range("I141").Value = 1248.5
$ws.Range("J141").Value = 2400
$ws.Range("K141").Value = 3745.5
$ws.Range("L141").Value = 7200
$ws.Range("M141").Value = 1434.5
$ws.Range("N141").Value = -17560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 48060.863
$ws.Range("I2").Value = 54386.26
$ws.Range("J2").Value = 8000
$ws.Range("K2").Value = 54386.26
$ws.Range("L2").Value = 8000
$ws.Range("M2").Value = -54273.26
$ws.Range("N2").Value = -8226

$ws.Range("H32").Value = 17033.842
$ws.Range("I32").Value = 14585
$ws.Range("J32").Value = 72745
$ws.Range("K32").Value = 14585
$ws.Range("L32").Value = 72745
$ws.Range("M32").Value = -14298
$ws.Range("N32").Value = -73319

$ws.Range("H61").Value = 1577.5714
$ws.Range("I61").Value = 1334.2174
$ws.Range("J61").Value = 2044
$ws.Range("K61").Value = 1334.2174
$ws.Range("L61").Value = 2044
$ws.Range("M61").Value = -1122.2174
$ws.Range("N61").Value = -2468

$ws.Range("H101").Value = 54900
$ws.Range("J101").Value = 54900
$ws.Range("L101").Value = 54900
$ws.Range("N101").Value = -61390

$ws.Range("H102").Value = 3700.3333
$ws.Range("I102").Value = 3700.3333
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3700.3333
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2078.3333
$ws.Range("N102").ClearContents()

$ws.Range("H104").Value = 29537.5
$ws.Range("J104").Value = 29537.5
$ws.Range("L104").Value = 29537.5
$ws.Range("N104").Value = -36525.5

$ws.Range("H110").Value = 666.64703
$ws.Range("I110").Value = 634.6429000000001
$ws.Range("J110").Value = 816
$ws.Range("K110").Value = 634.6429000000001
$ws.Range("L110").Value = 816
$ws.Range("M110").Value = 1410.3571
$ws.Range("N110").Value = -4906

$ws.Range("H116").Value = 48060.863
$ws.Range("I116").Value = 54386.26
$ws.Range("J116").Value = 8000
$ws.Range("K116").Value = 54386.26
$ws.Range("L116").Value = 8000
$ws.Range("M116").Value = -52092.26
$ws.Range("N116").Value = -12588

$ws.Range("H122").Value = 1093.1666
$ws.Range("I122").Value = 1010.7727
$ws.Range("J122").Value = 1999.5
$ws.Range("K122").Value = 3032.3181
$ws.Range("L122").Value = 5998.5
$ws.Range("M122").Value = -582.3181
$ws.Range("N122").Value = -10898.5

$ws.Range("H136").Value = 1577.5714
$ws.Range("I136").Value = 1334.2174
$ws.Range("J136").Value = 2044
$ws.Range("K136").Value = 4002.6522
$ws.Range("L136").Value = 6132
$ws.Range("M136").Value = -1452.6522
$ws.Range("N136").Value = -11232

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 48060.863
$ws.Range("I3").Value = 54386.26
$ws.Range("J3").Value = 8000
$ws.Range("K3").Value = 54386.26
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = -54272.26
$ws.Range("N3").Value = -8228

$ws.Range("H105").Value = 2974.6843
$ws.Range("I105").Value = 3065.8235
$ws.Range("J105").Value = 2200
$ws.Range("K105").Value = 3065.8235
$ws.Range("L105").Value = 2200
$ws.Range("M105").Value = -1318.8235
$ws.Range("N105").Value = -5694

$ws.Range("H107").Value = 1685.1333
$ws.Range("I107").Value = 1648.909
$ws.Range("J107").Value = 1784.75
$ws.Range("K107").Value = 1648.909
$ws.Range("L107").Value = 1784.75
$ws.Range("M107").Value = 271.0909999999999
$ws.Range("N107").Value = -5624.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4193.2
$ws.Range("I58").Value = 843.6
$ws.Range("J58").Value = 6705.4
$ws.Range("K58").Value = 843.6
$ws.Range("L58").Value = 6705.4
$ws.Range("M58").Value = -640.6
$ws.Range("N58").Value = -7111.4

$ws.Range("H122").Value = 527046.0600000001
$ws.Range("I122").Value = 1000779.1
$ws.Range("J122").Value = 676
$ws.Range("K122").Value = 3002337.3
$ws.Range("L122").Value = 2028
$ws.Range("M122").Value = -2999887.3
$ws.Range("N122").Value = -6928

$ws.Range("H134").Value = 2361.35
$ws.Range("I134").Value = 1620
$ws.Range("J134").Value = 3901.077
$ws.Range("K134").Value = 4860
$ws.Range("L134").Value = 11703.231
$ws.Range("M134").Value = -2325
$ws.Range("N134").Value = -16773.231

$ws.Range("H136").Value = 4193.2
$ws.Range("I136").Value = 843.6
$ws.Range("J136").Value = 6705.4
$ws.Range("K136").Value = 2530.8
$ws.Range("L136").Value = 20116.2
$ws.Range("M136").Value = 19.19999999999982
$ws.Range("N136").Value = -25216.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 54.857143
$ws.Range("I23").Value = 48.166668
$ws.Range("J23").Value = 59.875
$ws.Range("K23").Value = 144.500004
$ws.Range("L23").Value = 179.625
$ws.Range("M23").Value = 90.49999600000001
$ws.Range("N23").Value = -649.625

$ws.Range("H68").Value = 1749.6207
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1749.6207
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 5248.8621
$ws.Range("N68").Value = -6870.8621
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 1749.6207
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1749.6207
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 15746.5863
$ws.Range("N71").Value = -23858.5863
$ws.Range("M71").ClearContents()

$ws.Range("H107").Value = 695.1177
$ws.Range("I107").Value = 223.16129
$ws.Range("J107").Value = 1426.65
$ws.Range("K107").Value = 669.48387
$ws.Range("L107").Value = 4279.950000000001
$ws.Range("M107").Value = 1250.51613
$ws.Range("N107").Value = -8119.950000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4437.517
$ws.Range("I70").Value = 4080.3809
$ws.Range("J70").Value = 5375
$ws.Range("K70").Value = 4080.3809
$ws.Range("L70").Value = 5375
$ws.Range("M70").Value = -3810.3809
$ws.Range("N70").Value = -5915

$ws.Range("H73").Value = 4437.517
$ws.Range("I73").Value = 4080.3809
$ws.Range("J73").Value = 5375
$ws.Range("K73").Value = 4080.3809
$ws.Range("L73").Value = 5375
$ws.Range("M73").Value = -3144.3809
$ws.Range("N73").Value = -7247

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 967
$ws.Range("I61").Value = 967
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 967
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -765
$ws.Range("N61").ClearContents()

$ws.Range("H113").Value = 967
$ws.Range("I113").Value = 967
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 967
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1203
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 365.27585
$ws.Range("I107").Value = 318.5
$ws.Range("J107").Value = 512.2857
$ws.Range("K107").Value = 955.5
$ws.Range("L107").Value = 1536.8571
$ws.Range("M107").Value = 964.5
$ws.Range("N107").Value = -5376.8571

$ws.Range("H135").Value = 53088.176
$ws.Range("J135").Value = 53088.176
$ws.Range("L135").Value = 53088.176
$ws.Range("N135").Value = -63228.176
